# Auto-generated PowerShell COM-interop script.
# Applies the "Updated the test case document with test results" edit:
#  - Inserts w:proofErr (spellStart/spellEnd, gramStart/gramEnd) markers that
#    Word's proofing tools add, splitting the runs they wrap around.
#  - Fills in "Pass" results for five previously-empty result cells.
#  - Relocates the _GoBack bookmark from the last table cell to the first
#    empty paragraph following the third table.
$d = $word.ActiveDocument
$paras = $d.Paragraphs

$paras.Item(35).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="7B36FCE0" w14:textId="5CA67E5B" w:rsidR="00641E0B" w:rsidRDefault="00641E0B"><w:r><w:t xml:space="preserve">value: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>abc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$paras.Item(83).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="7203CDDA" w14:textId="55C75A51" w:rsidR="00222B77" w:rsidRDefault="00222B77" w:rsidP="00222B77"><w:r><w:t xml:space="preserve">To test if inputting </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">0 </w:t></w:r><w:r><w:t xml:space="preserve"> will</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> stop the python script</w:t></w:r></w:p>')
$paras.Item(103).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="12434080" w14:textId="77777777" w:rsidR="00354CBA" w:rsidRDefault="00354CBA" w:rsidP="00354CBA"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">“Enter the name of the data </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>file :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> ____”</w:t></w:r></w:p>')
$paras.Item(132).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="560FD099" w14:textId="77777777" w:rsidR="00354CBA" w:rsidRDefault="00354CBA" w:rsidP="00354CBA"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>Mazz.cev</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$paras.Item(134).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="3544A00C" w14:textId="77777777" w:rsidR="00354CBA" w:rsidRDefault="00354CBA" w:rsidP="00354CBA"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:lastRenderedPageBreak/><w:t>aofjb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')
$paras.Item(152).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="0ABA4C5A" w14:textId="3B05AA4F" w:rsidR="00354997" w:rsidRDefault="000B6F6B" w:rsidP="001F3885"><w:proofErr w:type="spellStart"/><w:r><w:t>emptyList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = []</w:t></w:r></w:p>')
$paras.Item(154).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="002150C8" w14:textId="77777777" w:rsidR="00354997" w:rsidRDefault="00354997" w:rsidP="001F3885"><w:r><w:t>Pass</w:t></w:r></w:p>')
$paras.Item(160).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="61F416FB" w14:textId="3507EC64" w:rsidR="00354997" w:rsidRDefault="000B6F6B" w:rsidP="001F3885"><w:proofErr w:type="spellStart"/><w:r><w:t>notmazeList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = ["1","2"]</w:t></w:r></w:p>')
$paras.Item(162).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="025D486E" w14:textId="77777777" w:rsidR="00354997" w:rsidRDefault="00354997" w:rsidP="001F3885"><w:r><w:t>Pass</w:t></w:r></w:p>')
$paras.Item(167).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="0719EE59" w14:textId="457A0442" w:rsidR="000B6F6B" w:rsidRDefault="000B6F6B" w:rsidP="000B6F6B"><w:proofErr w:type="spellStart"/><w:r><w:t>mazeList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = ["X","O","A","B"]</w:t></w:r></w:p>')
$paras.Item(169).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="1C480B85" w14:textId="77777777" w:rsidR="000B6F6B" w:rsidRDefault="000B6F6B" w:rsidP="000B6F6B"><w:r><w:t>Pass</w:t></w:r></w:p>')
$paras.Item(174).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="22C07D41" w14:textId="7298258A" w:rsidR="000B6F6B" w:rsidRDefault="000B6F6B" w:rsidP="000B6F6B"><w:proofErr w:type="spellStart"/><w:r><w:t>mazeList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = ["X","O","A","B"]</w:t></w:r></w:p>')
$paras.Item(176).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="6F9673F5" w14:textId="77777777" w:rsidR="000B6F6B" w:rsidRDefault="000B6F6B" w:rsidP="000B6F6B"><w:r><w:t>Pass</w:t></w:r></w:p>')
$paras.Item(181).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="3ED2CA73" w14:textId="5D422239" w:rsidR="000B6F6B" w:rsidRDefault="000B6F6B" w:rsidP="000B6F6B"><w:proofErr w:type="spellStart"/><w:r><w:t>notmazeList</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = ["1","2"]</w:t></w:r></w:p>')
$paras.Item(183).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="13546B46" w14:textId="77777777" w:rsidR="000B6F6B" w:rsidRDefault="000B6F6B" w:rsidP="000B6F6B"><w:r><w:t>Pass</w:t></w:r></w:p>')
$paras.Item(185).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="6D91AB1E" w14:textId="6CB56E48" w:rsidR="00354997" w:rsidRDefault="00354997"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
$paras.Item(285).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="21591577" w14:textId="0881738E" w:rsidR="00214E4B" w:rsidRDefault="00214E4B" w:rsidP="001F3885"><w:r><w:t>Test value: [</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>X,O</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>,O,A,O,O,X]</w:t></w:r></w:p>')
$paras.Item(292).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="37356860" w14:textId="737A1FFA" w:rsidR="00550F6A" w:rsidRDefault="00550F6A" w:rsidP="001F3885"><w:r><w:t xml:space="preserve">Test value: </w:t></w:r><w:r><w:t>[</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>X,</w:t></w:r><w:r><w:t>B</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>,O,A,O,O,X]</w:t></w:r></w:p>')
$paras.Item(343).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="3E08E574" w14:textId="77777777" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigOption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>')
$paras.Item(364).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="306A5338" w14:textId="77777777" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigOption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: 2</w:t></w:r></w:p>')
$paras.Item(372).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="06D2EA04" w14:textId="694D27E2" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:r><w:t xml:space="preserve">To test the creation of a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> point</w:t></w:r></w:p>')
$paras.Item(374).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="2C3289DB" w14:textId="77777777" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigOption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: 3</w:t></w:r></w:p>')
$paras.Item(403).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="330E9458" w14:textId="77777777" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigOption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: 1</w:t></w:r></w:p>')
$paras.Item(413).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="389FBF28" w14:textId="77777777" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigOption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>')
$paras.Item(424).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="17ED2F31" w14:textId="77777777" w:rsidR="00A75E1C" w:rsidRDefault="00A75E1C" w:rsidP="00A75E1C"><w:proofErr w:type="spellStart"/><w:r><w:t>ConfigOption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>')
$paras.Item(429).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="2E992D21" w14:textId="4E29A312" w:rsidR="00A75E1C" w:rsidRDefault="007B5772" w:rsidP="00A75E1C"><w:r><w:t>Pass</w:t></w:r></w:p>')
